$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42, shifting existing rows 42-98 down to 43-99.
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with the new weekly record.
$ws.Cells.Item(42, 1).Value = 7
$ws.Cells.Item(42, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(42, 3).Value = "Ñuble"
$ws.Cells.Item(42, 4).Value = 45117
$ws.Cells.Item(42, 5).Value = 16
$ws.Cells.Item(42, 6).Value = 100112013
$ws.Cells.Item(42, 7).Value = "Alcachofa"
$ws.Cells.Item(42, 8).Value = "Argentina(o)"
$ws.Cells.Item(42, 9).Value = "Primera"
$ws.Cells.Item(42, 10).Value = 50
$ws.Cells.Item(42, 11).Value = 17000
$ws.Cells.Item(42, 12).Value = 17000
$ws.Cells.Item(42, 13).Value = 17000
$ws.Cells.Item(42, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(42, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(42, 16).Value = 340
$ws.Cells.Item(42, 17).Value = 50
$ws.Cells.Item(42, 18).Value = "Hortaliza"
